$d = $word.ActiveDocument

# This cover-letter paragraph originally reads (as a single run):
#   "...compared with regard to certain research goals (e.g., in
#   psycholinguistics). The goal of the present study..."
#
# It needs to become (split across several runs with identical
# Times New Roman / 22pt formatting, matching the author's incremental
# edit history):
#   "...compared with regard to certain methodological goals (e.g.,
#   removal of \u201csocially-irrelevant variation\u201d in sociolinguistic
#   research). The goal of the present study..."
#
# Simple Find/Replace merges same-formatted text back into one run, so
# each edit below is done by briefly toggling Bold on/off around the
# Range.Text assignment -- this forces the engine to keep the edited
# slice as its own <w:r> instead of re-merging it with its identically
# formatted neighbours.

# --- Step 1: "research" -> "methodological" ------------------------------
$rng = $d.Content
$rng.Find.Execute("certain research goals", $true, $false, $false, $false, `
                   $false, $true, 1, $false, "", 0) | Out-Null
$r1 = $d.Range($rng.Start + 8, $rng.Start + 16)   # just the word "research"
$r1.Bold = 1
$r1.Text = "methodological"
$r1.Bold = 0

# --- Step 2: "in psycholinguistics" -> full replacement clause -----------
$rng2 = $d.Content
$rng2.Find.Execute("in psycholinguistics", $true, $false, $false, $false, `
                    $false, $true, 1, $false, "", 0) | Out-Null
$r2 = $d.Range($rng2.Start, $rng2.End)
$r2.Bold = 1
$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$newClause = "removal of " + $quoteOpen + "socially-irrelevant variation" + `
             $quoteClose + " in sociolinguistic research"
$r2.Text = $newClause
$r2.Bold = 0

# --- Step 3: split "in" back out into its own run -------------------------
$anchorText = "variation" + $quoteClose + " in sociolinguistic"
$rng3 = $d.Content
$rng3.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$prefixLen = ("variation" + $quoteClose + " ").Length
$inStart = $rng3.Start + $prefixLen
$inEnd = $inStart + 2
$r3 = $d.Range($inStart, $inEnd)
$r3.Bold = 1
$r3.Text = "in"
$r3.Bold = 0
